$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 92: 2020-02-27
$ws.Cells.Item(92, 1).Value = 1582761600
Set-TextCell $ws.Cells.Item(92, 2) "2020-02-27"
Set-TextCell $ws.Cells.Item(92, 3) "0212"
Set-TextCell $ws.Cells.Item(92, 4) "SDS"
$ws.Cells.Item(92, 5).Value = 0.2
$ws.Cells.Item(92, 6).Value = 0.22
$ws.Cells.Item(92, 7).Value = 0.19
$ws.Cells.Item(92, 8).Value = 0.22
$ws.Cells.Item(92, 9).Value = 1442000

# Row 93: 2020-02-28
$ws.Cells.Item(93, 1).Value = 1582848000
Set-TextCell $ws.Cells.Item(93, 2) "2020-02-28"
Set-TextCell $ws.Cells.Item(93, 3) "0212"
Set-TextCell $ws.Cells.Item(93, 4) "SDS"
$ws.Cells.Item(93, 5).Value = 0.21
$ws.Cells.Item(93, 6).Value = 0.21
$ws.Cells.Item(93, 7).Value = 0.195
$ws.Cells.Item(93, 8).Value = 0.2
$ws.Cells.Item(93, 9).Value = 370000
